$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -1
    4  = -5
    5  = -3
    6  = -4
    7  = 1
    8  = -6
    9  = -3
    10 = -3
    11 = -4
    12 = 1
    13 = 4
    14 = 4
    15 = -2
    16 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
